$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 18.47453410478584
$ws.Range("B2").Value = 12.15840038389621
$ws.Range("C2").Value = 23.73890881088179
$ws.Range("A3").Value = 20.69927875196798
$ws.Range("B3").Value = 16.12360934460579
$ws.Range("C3").Value = 24.68925092425341
$ws.Range("A4").Value = 20.31579241777433
$ws.Range("B4").Value = 15.65180638030963
$ws.Range("C4").Value = 23.91459789473231
$ws.Range("A5").Value = 25.19350774608662
$ws.Range("B5").Value = 20.01858622958437
$ws.Range("C5").Value = 29.7536787166161
$ws.Range("A6").Value = 24.47270220256332
$ws.Range("B6").Value = 18.51811415773333
$ws.Range("C6").Value = 29.71866506913102
$ws.Range("A7").Value = 16.74230273519213
$ws.Range("B7").Value = 11.57602230804374
$ws.Range("C7").Value = 21.5590412241435
$ws.Range("A8").Value = 9.613837263511396
$ws.Range("B8").Value = 5.334797848705129
$ws.Range("C8").Value = 14.09229993322819
$ws.Range("A9").Value = 22.53343943542593
$ws.Range("B9").Value = 17.34581389443653
$ws.Range("C9").Value = 27.29305688722882
$ws.Range("A10").Value = 9.985389128790441
$ws.Range("B10").Value = 4.862945197225781
$ws.Range("C10").Value = 14.93085509888761
$ws.Range("A11").Value = 15.78375775608151
$ws.Range("B11").Value = 9.059611242895151
$ws.Range("C11").Value = 21.69902176443171
$ws.Range("A12").Value = 23.5761827849423
$ws.Range("B12").Value = 15.67812248077889
$ws.Range("C12").Value = 31.1374495987938
$ws.Range("A13").Value = 12.12730704303809
$ws.Range("B13").Value = 5.022912108970725
$ws.Range("C13").Value = 18.60620167464742
$ws.Range("A14").Value = 13.33393133206784
$ws.Range("B14").Value = 7.093135776746295
$ws.Range("C14").Value = 19.53202552748217
$ws.Range("A15").Value = 28.70658342781184
$ws.Range("B15").Value = 21.90682976235371
$ws.Range("C15").Value = 34.8017527147783
$ws.Range("A16").Value = 24.43374594214981
$ws.Range("B16").Value = 18.37715526893616
$ws.Range("C16").Value = 29.58913328901266
$ws.Range("A17").Value = 23.66726532860142
$ws.Range("B17").Value = 16.26989691901897
$ws.Range("C17").Value = 30.47784804171852
$ws.Range("A18").Value = 20.29639955216287
$ws.Range("B18").Value = 15.64024316297008
$ws.Range("C18").Value = 23.90740259015658
$ws.Range("A19").Value = 16.40746201441004
$ws.Range("B19").Value = 8.687921097340558
$ws.Range("C19").Value = 24.23974881310412
$ws.Range("A20").Value = 26.94254793169214
$ws.Range("B20").Value = 19.33853268067864
$ws.Range("C20").Value = 34.07402041307092
$ws.Range("A21").Value = 27.50687703928825
$ws.Range("B21").Value = 21.37140808661711
$ws.Range("C21").Value = 33.05494165152641
$ws.Range("A22").Value = 10.60745770024984
$ws.Range("B22").Value = 5.045099259466177
$ws.Range("C22").Value = 16.27946326543118
$ws.Range("A23").Value = 30.27397652764357
$ws.Range("B23").Value = 23.15707734968162
$ws.Range("C23").Value = 36.83926127608205
$ws.Range("A24").Value = 24.43361137443287
$ws.Range("B24").Value = 18.41428501713848
$ws.Range("C24").Value = 29.57863333728151
$ws.Range("A25").Value = 20.90608493956147
$ws.Range("B25").Value = 15.57031038888938
$ws.Range("C25").Value = 25.27508449703469
$ws.Range("A26").Value = 27.6998310035032
$ws.Range("B26").Value = 21.32129850465838
$ws.Range("C26").Value = 33.65033986061898
$ws.Range("A27").Value = 12.78481961367045
$ws.Range("B27").Value = 4.980935237562744
$ws.Range("C27").Value = 20.78903088820732
$ws.Range("A28").Value = 17.70047528569038
$ws.Range("B28").Value = 12.17049151683925
$ws.Range("C28").Value = 22.93122711681796
$ws.Range("A29").Value = 32.20812374554371
$ws.Range("B29").Value = 24.03773961837921
$ws.Range("C29").Value = 39.43368261944684
$ws.Range("A30").Value = 22.55202366977467
$ws.Range("B30").Value = 15.96510352916992
$ws.Range("C30").Value = 28.32859392310145
$ws.Range("A31").Value = 13.05967691509111
$ws.Range("B31").Value = 6.93618352165032
$ws.Range("C31").Value = 18.57583385781128
$ws.Range("A32").Value = 10.9945607980956
$ws.Range("B32").Value = 6.500940091135796
$ws.Range("C32").Value = 15.8303504910484
$ws.Range("A33").Value = 13.47393129040054
$ws.Range("B33").Value = 7.338669279392564
$ws.Range("C33").Value = 19.29556175034782
$ws.Range("A34").Value = 19.81657141236616
$ws.Range("B34").Value = 11.93275845174982
$ws.Range("C34").Value = 27.47706320234663
$ws.Range("A35").Value = 25.39631057074351
$ws.Range("B35").Value = 18.76277557784498
$ws.Range("C35").Value = 30.92541926651926
$ws.Range("A36").Value = 16.73615157229322
$ws.Range("B36").Value = 11.70752942363443
$ws.Range("C36").Value = 21.38235966681745
$ws.Range("A37").Value = 29.31794840790199
$ws.Range("B37").Value = 21.86108562213212
$ws.Range("C37").Value = 35.94844393382063
$ws.Range("A38").Value = 14.49187937023411
$ws.Range("B38").Value = 7.533407113693881
$ws.Range("C38").Value = 21.39436623733113
$ws.Range("A39").Value = 24.17243399364451
$ws.Range("B39").Value = 17.10327535258527
$ws.Range("C39").Value = 31.79179330120244
$ws.Range("A40").Value = 24.21146517580116
$ws.Range("B40").Value = 18.66431322254551
$ws.Range("C40").Value = 28.87993533120542
$ws.Range("A41").Value = 15.66052117592052
$ws.Range("B41").Value = 11.18422630364136
$ws.Range("C41").Value = 20.09448568560348
$ws.Range("A42").Value = 16.09421480737154
$ws.Range("B42").Value = 9.780909348882759
$ws.Range("C42").Value = 22.21550553880288
$ws.Range("A43").Value = 22.5527866277756
$ws.Range("B43").Value = 15.94620221459947
$ws.Range("C43").Value = 28.30399307553023
$ws.Range("A44").Value = 13.60482624074503
$ws.Range("B44").Value = 9.001622333057654
$ws.Range("C44").Value = 17.76655863288588
$ws.Range("A45").Value = 16.95044057077782
$ws.Range("B45").Value = 10.44860351656015
$ws.Range("C45").Value = 24.02009461033913
$ws.Range("A46").Value = 12.31262839425655
$ws.Range("B46").Value = 7.234325810281641
$ws.Range("C46").Value = 17.04977206542621
$ws.Range("A47").Value = 10.88285966439648
$ws.Range("B47").Value = 4.555652730229136
$ws.Range("C47").Value = 17.16913693726036
$ws.Range("A48").Value = 23.6084153290437
$ws.Range("B48").Value = 15.7499856885759
$ws.Range("C48").Value = 31.19390603122809
$ws.Range("A49").Value = 22.04927863846541
$ws.Range("B49").Value = 17.12919413104799
$ws.Range("C49").Value = 26.36772188528159
$ws.Range("A50").Value = 29.12547256439923
$ws.Range("B50").Value = 21.51840229715238
$ws.Range("C50").Value = 35.84227765733107
$ws.Range("A51").Value = 20.71627673966479
$ws.Range("B51").Value = 16.12880032756508
$ws.Range("C51").Value = 24.70080677618954
$ws.Range("A52").Value = 11.8743807951967
$ws.Range("B52").Value = 7.304541432368743
$ws.Range("C52").Value = 16.10979117430914
$ws.Range("A53").Value = 24.22573719013922
$ws.Range("B53").Value = 17.35141147657535
$ws.Range("C53").Value = 31.64386707181844
$ws.Range("A54").Value = 28.84507553898488
$ws.Range("B54").Value = 21.9143139934278
$ws.Range("C54").Value = 34.98368190490864
$ws.Range("A55").Value = 27.74152265453415
$ws.Range("B55").Value = 20.95179187029396
$ws.Range("C55").Value = 33.31416531448477
$ws.Range("A56").Value = 13.01684693258103
$ws.Range("B56").Value = 7.203034533340197
$ws.Range("C56").Value = 18.64129468933451
$ws.Range("A57").Value = 10.9945607980956
$ws.Range("B57").Value = 6.500940091135796
$ws.Range("C57").Value = 15.8303504910484
$ws.Range("A58").Value = 29.27627171089414
$ws.Range("B58").Value = 22.31669803358166
$ws.Range("C58").Value = 35.63809765937265
$ws.Range("A59").Value = 14.77627736571667
$ws.Range("B59").Value = 9.028798426925272
$ws.Range("C59").Value = 20.32697394163779
$ws.Range("A60").Value = 30.02388602507967
$ws.Range("B60").Value = 22.64657327369743
$ws.Range("C60").Value = 36.82523547901404
$ws.Range("A61").Value = 13.82853365080115
$ws.Range("B61").Value = 7.647152705208994
$ws.Range("C61").Value = 19.64864491776657
$ws.Range("A62").Value = 15.80761741145067
$ws.Range("B62").Value = 9.228955496495228
$ws.Range("C62").Value = 21.64422216632413
$ws.Range("A63").Value = 14.44145019087484
$ws.Range("B63").Value = 7.730723531660597
$ws.Range("C63").Value = 22.1713894330229
$ws.Range("A64").Value = 30.06625723098194
$ws.Range("B64").Value = 22.65005458868843
$ws.Range("C64").Value = 36.8618267285462
$ws.Range("A65").Value = 10.49762901288996
$ws.Range("B65").Value = 6.398557959679606
$ws.Range("C65").Value = 15.03607425900588
$ws.Range("A66").Value = 19.5814779793238
$ws.Range("B66").Value = 11.71369304427932
$ws.Range("C66").Value = 27.07108568088774
$ws.Range("A67").Value = 19.25963022363516
$ws.Range("B67").Value = 12.63348529102075
$ws.Range("C67").Value = 25.10167968610517
$ws.Range("A68").Value = 25.82377520280117
$ws.Range("B68").Value = 19.80462145683063
$ws.Range("C68").Value = 31.48243684444612
$ws.Range("A69").Value = 11.24608654119153
$ws.Range("B69").Value = 4.786728753815812
$ws.Range("C69").Value = 17.18131032210431
$ws.Range("A70").Value = 24.13254753763742
$ws.Range("B70").Value = 18.73806096368129
$ws.Range("C70").Value = 29.09301285792912
$ws.Range("A71").Value = 17.65703593128963
$ws.Range("B71").Value = 13.26032951098136
$ws.Range("C71").Value = 21.76138442850016
$ws.Range("A72").Value = 20.29658012299004
$ws.Range("B72").Value = 15.07842494988894
$ws.Range("C72").Value = 25.19781970757816
